# Fix Education alignment:
#  - The paragraph that used to hold four plain tabs followed by the
#    "{-w:p educations}{allEducation}{/educations}" template tags is split
#    differently: the first two tabs now stay (bolded) at the end of the
#    "{educationsLabel}" paragraph, and the template-tag paragraph gets a
#    2880-twip left indent (matching the "{title}" paragraph) plus bold
#    paragraph mark formatting, so the templated rows line up correctly.
#  - The stray "_GoBack" bookmark that Word had left sitting in the grants
#    paragraph (right before "{/grants}") is relocated to the end of the
#    newly reshaped education paragraph.

$d = $word.ActiveDocument

# --- Remove the old _GoBack bookmark (sits just before "{/grants}") ---
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

# --- Locate the two paragraphs that make up the Education heading block ---
# Paragraph N   : "{educationsLabel} "
# Paragraph N+1 : four tabs + "{-w:p educations}{allEducation}{/educations}"
$eduLabelPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $txt = $d.Paragraphs($i).Range.Text
  if ($txt -like "*{educationsLabel}*") {
    $eduLabelPara = $i
    break
  }
}

$tabsPara = $eduLabelPara + 1

$p1 = $d.Paragraphs($eduLabelPara)
$p2 = $d.Paragraphs($tabsPara)

$rangeStart = $p1.Range.Start
$rangeEnd = $p2.Range.End
$r = $d.Range($rangeStart, $rangeEnd)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = ''
$newXml += '<w:p ' + $wNs + '>'
$newXml +=   '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>'
$newXml +=   '<w:r><w:rPr><w:b/></w:rPr><w:t>{educationsLabel}</w:t></w:r>'
$newXml +=   '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$newXml +=   '<w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r>'
$newXml +=   '<w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r>'
$newXml += '</w:p>'
$newXml += '<w:p ' + $wNs + '>'
$newXml +=   '<w:pPr><w:ind w:left="2880"/><w:rPr><w:b/></w:rPr></w:pPr>'
$newXml +=   '<w:r><w:t>{</w:t></w:r>'
$newXml +=   '<w:r><w:t>-w:p education</w:t></w:r>'
$newXml +=   '<w:r><w:t>s</w:t></w:r>'
$newXml +=   '<w:r><w:t>}</w:t></w:r>'
$newXml +=   '<w:r><w:t>{</w:t></w:r>'
$newXml +=   '<w:r><w:t>allEducation</w:t></w:r>'
$newXml +=   '<w:r><w:t>}{</w:t></w:r>'
$newXml +=   '<w:r><w:t>/</w:t></w:r>'
$newXml +=   '<w:r><w:t>education</w:t></w:r>'
$newXml +=   '<w:r><w:t>s</w:t></w:r>'
$newXml +=   '<w:r><w:t>}</w:t></w:r>'
$newXml +=   '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$newXml +=   '<w:bookmarkEnd w:id="0"/>'
$newXml += '</w:p>'

$r.InsertXML($newXml)
